# Season update through 1/17
$wb = $excel.ActiveWorkbook

# --- 1) "Games" sheet: append the completed game (previously the next
#        scheduled game on 45306 vs SAS) as new row 40 with full box-score
#        stats. ---
$games = $wb.Worksheets.Item("Games")

$lastRow = 39
$newRow = $lastRow + 1

$games.Cells.Item($newRow, 1).Value  = 39        # Game
$games.Cells.Item($newRow, 2).Value  = 45306     # Date
$games.Cells.Item($newRow, 3).Value  = 1         # Streak
$games.Cells.Item($newRow, 4).Value  = 109       # Pts
$games.Cells.Item($newRow, 5).Value  = 105.1     # Pace
$games.Cells.Item($newRow, 6).Value  = 0.5       # eFG
$games.Cells.Item($newRow, 7).Value  = 13.5      # TOV
$games.Cells.Item($newRow, 8).Value  = 27.1      # ORB
$games.Cells.Item($newRow, 9).Value  = 0.172     # FTR
$games.Cells.Item($newRow, 10).Value = 103.7     # ORT
$games.Cells.Item($newRow, 11).Value = "SAS"     # OppID
$games.Cells.Item($newRow, 12).Value = 99        # OppPts
$games.Cells.Item($newRow, 13).Value = 0.458     # OppeFG
$games.Cells.Item($newRow, 14).Value = 14.2      # OppTOV
$games.Cells.Item($newRow, 15).Value = 17.3      # OppORB
$games.Cells.Item($newRow, 16).Value = 0.126     # OppFTR
$games.Cells.Item($newRow, 17).Value = 94.2      # OppORT
$games.Cells.Item($newRow, 18).Value = 1         # Location
$games.Cells.Item($newRow, 19).Value = 1         # Target

# Match the date-formatted style used by the rest of column B
$games.Cells.Item($newRow, 2).NumberFormat = $games.Cells.Item($lastRow, 2).NumberFormat

# --- 2) "Next" sheet: the game that was played (45306 vs SAS) is no
#        longer upcoming, so remove its row and shift the remaining
#        schedule up. ---
$next = $wb.Worksheets.Item("Next")
$next.Rows(2).Delete()
